# Add "URPbIC" (Union Representation Percentage by ISIC Code) row to the
# io-model section of the "Key to Variables" sheet, and tidy up the wording
# of the BGDP "optional" note so it no longer references "(in Vensim)".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")

# 1) Update the BGDP "optional" note in column G (row 140) to drop "(in Vensim) ".
$ws.Range("G140").Value = "You plan on displaying or using the emissions-per-unit-GDP or energy-use-per-unit-GDP graphs"

# 2) Insert a new row before the current row 149 (VAbIC), which will push
#    VAbIC (and everything below it) down by one row.
$ws.Rows.Item(149).Insert()

# 3) Populate the newly-inserted row 149 with the URPbIC entry.
$ws.Range("A149").Value = "io-model"
$ws.Range("B149").Value = "URPbIC"
$ws.Range("C149").Value = "Union Representation Percentage by ISIC Code"
$ws.Range("F149").Value = "medium"

# 4) The row-insert copies the format of the row above (TLIM, "high"/yellow-red
#    banding) cell-by-cell, so F149 comes out styled "high" instead of
#    "medium". Re-stamp F149's format from another "medium" cell in the same
#    banded section (F146, LPGRbIC) so the fill color matches the "medium"
#    importance rating.
$ws.Range("F146").Copy()
$ws.Range("F149").PasteSpecial(-4122)
$excel.CutCopyMode = $false
